$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1 -- copy G1's formatting (bold header style) to H1,
# then overwrite the value with "Save".
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data value for the Save column, row 2.
$ws.Range("H2").Value = 0
